$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.062.89"
$ws.Range("E2").Value = "  -0.83%  "
$ws.Range("D3").Value = "2.010.70"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("A1").Value = "'226.43"
$ws.Range("A1").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E5").Value = "  -1.95%  "
$ws.Range("A1").Value = "'0.607"
$ws.Range("A1").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("A1").Value = "'54.78"
$ws.Range("A1").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E8").Value = "  -3.92%  "
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("A1").Value = "'0.0787"
$ws.Range("A1").Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E10").Value = "  +2.10%  "
$ws.Range("D12").Value = "2.306.87"
$ws.Range("E12").Value = "  -1.90%  "
$ws.Range("A1").Value = "'14.16"
$ws.Range("A1").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E13").Value = "  -3.02%  "
$ws.Range("A1").Value = "'20.29"
$ws.Range("A1").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E14").Value = "  -1.77%  "
$ws.Range("A1").Value = "'0.741"
$ws.Range("A1").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E15").Value = "  -1.93%  "
$ws.Range("A1").Value = "'5.13"
$ws.Range("A1").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("D17").Value = "2.009.99"
$ws.Range("E17").Value = "  -1.55%  "
$ws.Range("D18").Value = "36.986.66"
$ws.Range("E18").Value = "  -0.94%  "
$ws.Range("A1").Value = "'6.11"
$ws.Range("A1").Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E19").Value = "  +0.24%  "
$ws.Range("A1").Value = "'68.82"
$ws.Range("A1").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("D21").Value = "0.0₃0818"
$ws.Range("E21").Value = "  -0.57%  "
$ws.Range("A1").Value = "'223.43"
$ws.Range("A1").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  +1.24%  "
$ws.Range("A1").Value = "'2.19"
$ws.Range("A1").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E25").Value = "  -5.91%  "
$ws.Range("A1").Value = "'165.95"
$ws.Range("A1").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E26").Value = "  -2.41%  "
$ws.Range("A1").Value = "'9.16"
$ws.Range("A1").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E27").Value = "  -7.18%  "
$ws.Range("A1").Value = "'18.69"
$ws.Range("A1").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E28").Value = "  -2.54%  "
$ws.Range("E29").Value = "  -3.69%  "
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("A1").Value = "'0.117"
$ws.Range("A1").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("A1").Value = "'4.51"
$ws.Range("A1").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E32").Value = "  -0.29%  "
$ws.Range("A1").Value = "'0.0610"
$ws.Range("A1").Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E33").Value = "  -1.76%  "
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("E35").Value = "  -5.62%  "
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("A1").Value = "'3.15"
$ws.Range("A1").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("A1").Value = "'5.34"
$ws.Range("A1").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("D40").Value = "1.478.25"
$ws.Range("E40").Value = "  +0.10%  "
$ws.Range("E41").Value = "  -3.36%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("A1").Value = "'16.55"
$ws.Range("A1").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("A1").Value = "'94.83"
$ws.Range("A1").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E43").Value = "  -3.41%  "
$ws.Range("A1").Value = "'0.0922"
$ws.Range("A1").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E44").Value = "  -3.26%  "
$ws.Range("E45").Value = "  -4.88%  "
$ws.Range("A1").Value = "'1.13"
$ws.Range("A1").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E46").Value = "  -4.53%  "
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("A1").Value = "'7.18"
$ws.Range("A1").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("E49").Value = "  -0.86%  "
$ws.Range("D50").Value = "2.195.25"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("A1").Value = "'44.41"
$ws.Range("A1").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("A1").Clear()
$ws.Range("E51").Value = "  -2.94%  "
